$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 406.76923
$ws.Range("I4").Value = 250
$ws.Range("K4").Value = 250
$ws.Range("M4").Value = -136
$ws.Range("H28").Value = 1824.75
$ws.Range("I28").Value = 849.5
$ws.Range("J28").Value = 2800
$ws.Range("K28").Value = 849.5
$ws.Range("L28").Value = 2800
$ws.Range("M28").Value = -364.5
$ws.Range("N28").Value = -3770
$ws.Range("H32").Value = 1185.875
$ws.Range("I32").Value = 1330
$ws.Range("J32").Value = 1099.4
$ws.Range("K32").Value = 1330
$ws.Range("L32").Value = 1099.4
$ws.Range("M32").Value = -1004
$ws.Range("N32").Value = -1751.4
$ws.Range("H33").Value = 415.66666
$ws.Range("I33").Value = 415.66666
$ws.Range("K33").Value = 415.66666
$ws.Range("M33").Value = -186.66666
$ws.Range("H58").Value = 924.8333
$ws.Range("J58").Value = 1750
$ws.Range("L58").Value = 5250
$ws.Range("N58").Value = -5550
$ws.Range("H92").Value = 431.75
$ws.Range("I92").Value = 487.7
$ws.Range("K92").Value = 487.7
$ws.Range("M92").Value = 760.3
$ws.Range("H98").Value = 1058.375
$ws.Range("I98").Value = 1155.9
$ws.Range("J98").Value = 895.8333
$ws.Range("K98").Value = 1155.9
$ws.Range("L98").Value = 895.8333
$ws.Range("M98").Value = 342.0999999999999
$ws.Range("N98").Value = -3891.8333
$ws.Range("H99").Value = 1271.8572
$ws.Range("I99").Value = 1345.6666
$ws.Range("K99").Value = 4036.9998
$ws.Range("M99").Value = -2538.9998
$ws.Range("H107").Value = 6300
$ws.Range("I107").Value = 6300
$ws.Range("K107").Value = 6300
$ws.Range("M107").Value = -4380
$ws.Range("H112").Value = 2206.8572
$ws.Range("J112").Value = 2291.5
$ws.Range("L112").Value = 6874.5
$ws.Range("N112").Value = -9090.5
$ws.Range("H122").Value = 1058.375
$ws.Range("I122").Value = 1155.9
$ws.Range("J122").Value = 895.8333
$ws.Range("K122").Value = 3467.7
$ws.Range("L122").Value = 2687.4999
$ws.Range("M122").Value = -1017.7
$ws.Range("N122").Value = -7587.4999
$ws.Range("H132").Value = 3068.9736
$ws.Range("I132").Value = 2013.5416
$ws.Range("K132").Value = 6040.6248
$ws.Range("M132").Value = -3510.6248
$ws.Range("H138").Value = 4075.077
$ws.Range("I138").Value = 3625.5715
$ws.Range("J138").Value = 4599.5
$ws.Range("K138").Value = 10876.7145
$ws.Range("L138").Value = 13798.5
$ws.Range("M138").Value = -5736.7145
$ws.Range("N138").Value = -24078.5
$ws.Range("H141").Value = 2524.625
$ws.Range("I141").Value = 2570.2856
$ws.Range("K141").Value = 7710.8568
$ws.Range("M141").Value = -2530.8568

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1106.3846
$ws.Range("J2").Value = 3525.6667
$ws.Range("L2").Value = 3525.6667
$ws.Range("N2").Value = -3751.6667
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = 0
$ws.Range("H102").Value = 2659
$ws.Range("I102").Value = 1432.3334
$ws.Range("K102").Value = 1432.3334
$ws.Range("M102").Value = 189.6666
$ws.Range("H116").Value = 1106.3846
$ws.Range("J116").Value = 3525.6667
$ws.Range("L116").Value = 3525.6667
$ws.Range("N116").Value = -8113.6667

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1106.3846
$ws.Range("J3").Value = 3525.6667
$ws.Range("L3").Value = 3525.6667
$ws.Range("N3").Value = -3753.6667
$ws.Range("H141").Value = 42648.5
$ws.Range("J141").Value = 42648.5
$ws.Range("L141").Value = 42648.5
$ws.Range("N141").Value = -53008.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1924.9474
$ws.Range("I7").Value = 1068.4166
$ws.Range("K7").Value = 1068.4166
$ws.Range("M7").Value = -955.4166
$ws.Range("H32").Value = 992.5
$ws.Range("I32").Value = 490
$ws.Range("K32").Value = 490
$ws.Range("M32").Value = -174
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").ClearContents()
$ws.Range("N100").Value = 0

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 45.739132
$ws.Range("I2").Value = 25
$ws.Range("J2").Value = 93.14286
$ws.Range("K2").Value = 150
$ws.Range("L2").Value = 558.85716
$ws.Range("M2").Value = -37
$ws.Range("N2").Value = -784.85716
$ws.Range("H32").Value = 487.5
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 487.5
$ws.Range("K32").Value = 0
$ws.Range("L32").ClearContents()
$ws.Range("M32").Value = 1462.5
$ws.Range("N32").Value = -2028.5
$ws.Range("H87").Value = 55
$ws.Range("I87").Value = 55
$ws.Range("K87").Value = 165
$ws.Range("M87").Value = 1083
$ws.Range("H90").Value = 55
$ws.Range("I90").Value = 55
$ws.Range("K90").Value = 495
$ws.Range("M90").Value = 5745
$ws.Range("H92").Value = 837.5
$ws.Range("I92").Value = 600
$ws.Range("K92").Value = 1800
$ws.Range("M92").Value = -552
$ws.Range("H109").Value = 3482.3333
$ws.Range("I109").Value = 223.5
$ws.Range("J109").Value = 10000
$ws.Range("K109").Value = 670.5
$ws.Range("L109").Value = 30000
$ws.Range("M109").Value = 369.5
$ws.Range("N109").Value = -32080
$ws.Range("H131").Value = 1159.9333
$ws.Range("J131").Value = 1164.2142
$ws.Range("L131").Value = 3492.6426
$ws.Range("N131").Value = -13572.6426

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 16724.455
$ws.Range("J43").Value = 26242.25
$ws.Range("L43").Value = 26242.25
$ws.Range("N43").Value = -26544.25
$ws.Range("H57").Value = 19998
$ws.Range("J57").Value = 19998
$ws.Range("L57").Value = 19998
$ws.Range("N57").Value = -21638
$ws.Range("H80").Value = 1199.5
$ws.Range("I80").Value = 1199.5
$ws.Range("K80").Value = 1199.5
$ws.Range("M80").Value = -201.5
$ws.Range("H83").Value = 1199.5
$ws.Range("I83").Value = 1199.5
$ws.Range("K83").Value = 5997.5
$ws.Range("M83").Value = -1005.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1814
$ws.Range("I46").Value = 1472.5
$ws.Range("J46").Value = 2360.4
$ws.Range("K46").Value = 1472.5
$ws.Range("L46").Value = 2360.4
$ws.Range("M46").Value = -1284.5
$ws.Range("N46").Value = -2736.4
$ws.Range("H100").Value = 7498.5
$ws.Range("I100").Value = 4997
$ws.Range("K100").Value = 4997
$ws.Range("M100").Value = -4456

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 15000
$ws.Range("J97").Value = 15000
$ws.Range("L97").Value = 15000
$ws.Range("N97").Value = -16982
